$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 142863040
$ws.Range("I18").Value = 5317.5
$ws.Range("J18").Value = 333340000
$ws.Range("K18").Value = 5317.5
$ws.Range("L18").Value = 333340000
$ws.Range("M18").Value = -5033.5
$ws.Range("N18").Value = -333340568
$ws.Range("H32").Value = 599
$ws.Range("I32").Value = 599
$ws.Range("K32").Value = 599
$ws.Range("M32").Value = -273
$ws.Range("H33").Value = 195.66667
$ws.Range("I33").Value = 181.3125
$ws.Range("K33").Value = 181.3125
$ws.Range("M33").Value = 47.6875
$ws.Range("H74").Value = 16667.54
$ws.Range("I74").Value = 16667.54
$ws.Range("K74").Value = 16667.54
$ws.Range("M74").Value = -15731.54
$ws.Range("H77").Value = 16667.54
$ws.Range("I77").Value = 16667.54
$ws.Range("K77").Value = 83337.70000000001
$ws.Range("M77").Value = -78657.70000000001
$ws.Range("H132").Value = 3242.9246
$ws.Range("I132").Value = 3366.9795
$ws.Range("K132").Value = 10100.9385
$ws.Range("M132").Value = -7570.9385
$ws.Range("H138").Value = 484614.72
$ws.Range("J138").Value = 738474.25
$ws.Range("L138").Value = 2215422.75
$ws.Range("N138").Value = -2225702.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H32").Value = 1471.7959
$ws.Range("I32").Value = 1445.7423
$ws.Range("K32").Value = 1445.7423
$ws.Range("M32").Value = -1158.7423
$ws.Range("H61").Value = 3012.4644
$ws.Range("I61").Value = 1771.75
$ws.Range("J61").Value = 6114.25
$ws.Range("K61").Value = 1771.75
$ws.Range("L61").Value = 6114.25
$ws.Range("M61").Value = -1559.75
$ws.Range("N61").Value = -6538.25
$ws.Range("H122").Value = 4494.037
$ws.Range("I122").Value = 3603.818
$ws.Range("J122").Value = 8411
$ws.Range("K122").Value = 10811.454
$ws.Range("L122").Value = 25233
$ws.Range("M122").Value = -8361.454000000002
$ws.Range("N122").Value = -30133
$ws.Range("H132").Value = 2116.42
$ws.Range("I132").Value = 1766.579
$ws.Range("J132").Value = 3224.25
$ws.Range("K132").Value = 5299.737
$ws.Range("L132").Value = 9672.75
$ws.Range("M132").Value = -2769.737
$ws.Range("N132").Value = -14732.75
$ws.Range("H136").Value = 3012.4644
$ws.Range("I136").Value = 1771.75
$ws.Range("J136").Value = 6114.25
$ws.Range("K136").Value = 5315.25
$ws.Range("L136").Value = 18342.75
$ws.Range("M136").Value = -2765.25
$ws.Range("N136").Value = -23442.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 5192.25
$ws.Range("I12").Value = 381.5
$ws.Range("J12").Value = 10003
$ws.Range("K12").Value = 381.5
$ws.Range("L12").Value = 10003
$ws.Range("M12").Value = -213.5
$ws.Range("N12").Value = -10339
$ws.Range("H86").Value = 6552.857
$ws.Range("I86").Value = 6842.75
$ws.Range("J86").Value = 6166.3335
$ws.Range("K86").Value = 6842.75
$ws.Range("L86").Value = 6166.3335
$ws.Range("M86").Value = -5719.75
$ws.Range("N86").Value = -8412.333500000001
$ws.Range("H89").Value = 6552.857
$ws.Range("I89").Value = 6842.75
$ws.Range("J89").Value = 6166.3335
$ws.Range("K89").Value = 34213.75
$ws.Range("L89").Value = 30831.6675
$ws.Range("M89").Value = -28597.75
$ws.Range("N89").Value = -42063.6675
$ws.Range("H107").Value = 1749.9445
$ws.Range("I107").Value = 1415.5385
$ws.Range("J107").Value = 2619.4
$ws.Range("K107").Value = 1415.5385
$ws.Range("L107").Value = 2619.4
$ws.Range("M107").Value = 504.4614999999999
$ws.Range("N107").Value = -6459.4
$ws.Range("H134").Value = 2857.9656
$ws.Range("I134").Value = 2146.6
$ws.Range("K134").Value = 6439.799999999999
$ws.Range("M134").Value = -3904.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I7").Value = 3929.7334
$ws.Range("J7").Value = 166667230
$ws.Range("K7").Value = 3929.7334
$ws.Range("L7").Value = 166667230
$ws.Range("M7").Value = -3816.7334
$ws.Range("N7").Value = -166667456
$ws.Range("H23").Value = 11666.667
$ws.Range("J23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15480
$ws.Range("H27").Value = 11666.667
$ws.Range("J27").Value = 15000
$ws.Range("L27").Value = 15000
$ws.Range("N27").Value = -15384
$ws.Range("H31").Value = 3601.9854
$ws.Range("I31").Value = 2890.98
$ws.Range("J31").Value = 5577
$ws.Range("K31").Value = 2890.98
$ws.Range("L31").Value = 5577
$ws.Range("M31").Value = -2595.98
$ws.Range("N31").Value = -6167
$ws.Range("H34").Value = 3601.9854
$ws.Range("I34").Value = 2890.98
$ws.Range("J34").Value = 5577
$ws.Range("K34").Value = 2890.98
$ws.Range("L34").Value = 5577
$ws.Range("M34").Value = -2688.98
$ws.Range("N34").Value = -5981
$ws.Range("H58").Value = 3308.5417
$ws.Range("I58").Value = 1922.4615
$ws.Range("J58").Value = 4946.636
$ws.Range("K58").Value = 1922.4615
$ws.Range("L58").Value = 4946.636
$ws.Range("M58").Value = -1719.4615
$ws.Range("N58").Value = -5352.636
$ws.Range("H86").Value = 2685.75
$ws.Range("I86").Value = 2223.2
$ws.Range("J86").Value = 4998.5
$ws.Range("K86").Value = 2223.2
$ws.Range("L86").Value = 4998.5
$ws.Range("M86").Value = -1100.2
$ws.Range("N86").Value = -7244.5
$ws.Range("H89").Value = 2685.75
$ws.Range("I89").Value = 2223.2
$ws.Range("J89").Value = 4998.5
$ws.Range("K89").Value = 11116
$ws.Range("L89").Value = 24992.5
$ws.Range("M89").Value = -5500
$ws.Range("N89").Value = -36224.5
$ws.Range("H107").Value = 702.5599999999999
$ws.Range("I107").Value = 705.4761999999999
$ws.Range("J107").Value = 687.25
$ws.Range("K107").Value = 705.4761999999999
$ws.Range("L107").Value = 687.25
$ws.Range("M107").Value = 1214.5238
$ws.Range("N107").Value = -4527.25
$ws.Range("H122").Value = 3086.6086
$ws.Range("I122").Value = 3025.8948
$ws.Range("J122").Value = 3375
$ws.Range("K122").Value = 9077.6844
$ws.Range("L122").Value = 10125
$ws.Range("M122").Value = -6627.6844
$ws.Range("N122").Value = -15025
$ws.Range("H132").Value = 4057.1
$ws.Range("I132").Value = 4164.579
$ws.Range("J132").Value = 3871.4546
$ws.Range("K132").Value = 12493.737
$ws.Range("L132").Value = 11614.3638
$ws.Range("M132").Value = -9963.736999999999
$ws.Range("N132").Value = -16674.3638
$ws.Range("H136").Value = 3308.5417
$ws.Range("I136").Value = 1922.4615
$ws.Range("J136").Value = 4946.636
$ws.Range("K136").Value = 5767.3845
$ws.Range("L136").Value = 14839.908
$ws.Range("M136").Value = -3217.3845
$ws.Range("N136").Value = -19939.908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 387.33334
$ws.Range("I14").Value = 387.33334
$ws.Range("K14").Value = 1162.00002
$ws.Range("M14").Value = -989.0000199999999
$ws.Range("H97").Value = 743184.1
$ws.Range("I97").Value = 1250249.2
$ws.Range("J97").Value = 67097.336
$ws.Range("K97").Value = 3750747.6
$ws.Range("L97").Value = 201292.008
$ws.Range("M97").Value = -3750251.6
$ws.Range("N97").Value = -202284.008
$ws.Range("H122").Value = 2419
$ws.Range("I122").Value = 2004
$ws.Range("J122").Value = 2470.875
$ws.Range("K122").Value = 18036
$ws.Range("L122").Value = 22237.875
$ws.Range("M122").Value = -15586
$ws.Range("N122").Value = -27137.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 525.5
$ws.Range("I97").Value = 586.7143
$ws.Range("K97").Value = 586.7143
$ws.Range("M97").Value = -90.71429999999998
$ws.Range("H122").Value = 1863.3125
$ws.Range("I122").Value = 1841
$ws.Range("J122").Value = 2198
$ws.Range("K122").Value = 5523
$ws.Range("L122").Value = 6594
$ws.Range("M122").Value = -3073
$ws.Range("N122").Value = -11494
$ws.Range("H132").Value = 3523.7046
$ws.Range("J132").Value = 5311.1
$ws.Range("L132").Value = 15933.3
$ws.Range("N132").Value = -20993.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8388.643
$ws.Range("I7").Value = 6320.5
$ws.Range("J7").Value = 9939.75
$ws.Range("K7").Value = 6320.5
$ws.Range("L7").Value = 9939.75
$ws.Range("M7").Value = -6208.5
$ws.Range("N7").Value = -10163.75
$ws.Range("H40").Value = 22286.701
$ws.Range("I40").Value = 37772.6
$ws.Range("J40").Value = 9730.566999999999
$ws.Range("K40").Value = 37772.6
$ws.Range("L40").Value = 9730.566999999999
$ws.Range("M40").Value = -37636.6
$ws.Range("N40").Value = -10002.567
$ws.Range("H82").Value = 5909.231
$ws.Range("J82").Value = 6996.6665
$ws.Range("L82").Value = 6996.6665
$ws.Range("N82").Value = -7718.6665
$ws.Range("H85").Value = 5909.231
$ws.Range("J85").Value = 6996.6665
$ws.Range("L85").Value = 6996.6665
$ws.Range("N85").Value = -9492.666499999999
$ws.Range("H122").Value = 5568.357
$ws.Range("I122").Value = 6029.143
$ws.Range("J122").Value = 5107.5713
$ws.Range("K122").Value = 18087.429
$ws.Range("L122").Value = 15322.7139
$ws.Range("M122").Value = -15637.429
$ws.Range("N122").Value = -20222.7139
$ws.Range("H126").Value = 8388.643
$ws.Range("I126").Value = 6320.5
$ws.Range("J126").Value = 9939.75
$ws.Range("K126").Value = 18961.5
$ws.Range("L126").Value = 29819.25
$ws.Range("M126").Value = -16491.5
$ws.Range("N126").Value = -34759.25
$ws.Range("H136").Value = 6191.4116
$ws.Range("I136").Value = 6225.4546
$ws.Range("J136").Value = 6129
$ws.Range("K136").Value = 18676.3638
$ws.Range("L136").Value = 18387
$ws.Range("M136").Value = -16126.3638
$ws.Range("N136").Value = -23487

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1706.037
$ws.Range("I126").Value = 1610.6957
$ws.Range("J126").Value = 2254.25
$ws.Range("K126").Value = 4832.0871
$ws.Range("L126").Value = 6762.75
$ws.Range("M126").Value = -2362.0871
$ws.Range("N126").Value = -11702.75

Write-Host "Applied all Gilgamesh Profits updates"